$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 24 de Septiembre de 2020 a las 02:15'
$ws.Range("B4").Value = 7138708
$ws.Range("C4").Value = 40771
$ws.Range("D4").Value = 4387230
$ws.Range("E4").Value = 2544920
$ws.Range("G4").Value = 1077
$ws.Range("H4").Value = 206558
$ws.Range("B6").Value = 4627780
$ws.Range("C6").Value = 32445
$ws.Range("E6").Value = 495829
$ws.Range("G6").Value = 906
$ws.Range("H6").Value = 139065
$ws.Range("B9").Value = 782695
$ws.Range("C9").Value = 6149
$ws.Range("D9").Value = 636489
$ws.Range("E9").Value = 114336
$ws.Range("G9").Value = 98
$ws.Range("H9").Value = 31870
$ws.Range("B13").Value = 664799
$ws.Range("C13").Value = 12625
$ws.Range("D13").Value = 525486
$ws.Range("E13").Value = 124937
$ws.Range("G13").Value = 424
$ws.Range("H13").Value = 14376
$ws.Range("D25").Value = 249500
$ws.Range("E25").Value = 20197
$ws.Range("A35").Value = 'Panama'
$ws.Range("B35").Value = 107990
$ws.Range("C35").Value = 706
$ws.Range("D35").Value = 84437
$ws.Range("E35").Value = 21262
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 2291
$ws.Range("A36").Value = 'Marruecos'
$ws.Range("B36").Value = 107743
$ws.Range("C36").Value = 2397
$ws.Range("D36").Value = 88244
$ws.Range("E36").Value = 17581
$ws.Range("G36").Value = 29
$ws.Range("H36").Value = 1918
$ws.Range("A37").Value = 'Kazajistan'
$ws.Range("B37").Value = 107450
$ws.Range("C37").Value = 76
$ws.Range("D37").Value = 102064
$ws.Range("E37").Value = 3687
$ws.Range("H37").Value = 1699
$ws.Range("B59").Value = 55464
$ws.Range("C59").Value = 2306
$ws.Range("E59").Value = 28200
$ws.Range("G59").Value = 24
$ws.Range("H59").Value = 555
$ws.Range("B127").Value = 4779
$ws.Range("C127").Value = 20
$ws.Range("D127").Value = 4560
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 101
$ws.Range("A150").Value = 'Benin'
$ws.Range("B150").Value = 2325
$ws.Range("C150").Value = 31
$ws.Range("D150").Value = 1954
$ws.Range("E150").Value = 331
$ws.Range("H150").Value = 40
$ws.Range("A151").Value = 'Guinea-Bisau'
$ws.Range("B151").Value = 2324
$ws.Range("D151").Value = 1549
$ws.Range("E151").Value = 736
$ws.Range("H151").Value = 39
$ws.Range("B154").Value = 1946
$ws.Range("C154").Value = 12
$ws.Range("D154").Value = 1661
$ws.Range("E154").Value = 238
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 47
$ws.Range("B155").Value = 1929
$ws.Range("C155").Value = 22
$ws.Range("E155").Value = 621
$ws.Range("B160").Value = 1654
$ws.Range("C160").Value = 36
$ws.Range("E160").Value = 263
$ws.Range("A214").Value = 'Islas Malvinas'
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
$ws.Range("A215").Value = 'Montserrat'
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
